# Add a new contact row (S.no=2, Name="ING Gabriela", Contact="529611701291")
# to the WhatsApp contact list, matching the formatting already used by row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply center alignment first so the new cells pick up the same plain
# "center, no border" style used elsewhere in the sheet.
$ws.Range("A3:C3").HorizontalAlignment = -4108

# Column C holds phone numbers stored as text (same as C2), so force the
# number format to Text before writing the value to avoid Excel coercing
# the long digit string into a numeric value.
$ws.Range("C3").NumberFormat = "@"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "ING Gabriela"
$ws.Range("C3").Value = "529611701291"

# Move the active selection, matching the saved workbook state.
$ws.Range("G7").Select()
